# Updated sensitivity-matrix values for Sheet1!B2:AF8 (7 rows x 31 cols)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(-2.978358019635732, 1.46110522798494, 1.76488022812725, -0.1819887215705478, -0.0649079407850771, -0.0007307741227984005, -7.557823177602677, 4.290037286492535, 5.37823968090039, -1.613803115652376, -0.4906743893445039, -0.005976284788717497, -4.882783070365717, 9.075332653468799, 12.54397119453234, -12.79800240606736, -3.829682313083997, -0.1088360584230113, 7.893833719290445, 9.347706739494658, 10.48938412452083, -20.01345706802031, -7.103493917402907, -0.6139735977514333, -18.65772957169785, 10.22202896656103, 6.459959547148602, 3.699109108859769, -0.8523620783265802, -0.8710059723279313, 0.2324312823225227),
    @(-2.978358019637952, 1.46110522798605, 1.76488022812836, -0.1819887215705132, -0.06490794078507493, -0.0007307741227983665, -7.557823177607116, 4.290037286496975, 5.37823968090483, -1.613803115652099, -0.4906743893444865, -0.005976284788717225, -4.882783070383476, 9.075332653477679, 12.5439711945501, -12.79800240606514, -3.829682313083441, -0.108836058422994, 7.89383371925492, 9.34770673951242, 10.48938412452971, -20.01345706801143, -7.103493917400686, -0.6139735977512946, -18.65772957173337, 10.22202896656547, 6.459959547157483, 3.699109108877533, -0.8523620783176982, -0.871005972326821, 0.2324312823226364),
    @(-3.225733044734957, 1.412074164503649, 1.646627961855945, 0.141646390334282, 0.02510842766352472, 0.0002761003779703362, -8.722718349465147, 3.609298674478238, 4.129622177110615, 0.7966373946401536, 0.1849327619165354, 0.00222734131918904, -16.11529125383388, 4.808850115160541, 4.989208589553469, 4.902080449535722, 1.376955804706298, 0.038196294851406, -9.456536497396705, 0.8474596023921247, 1.300286253585226, 5.250186914725215, 1.879489640070679, 0.1791140865794251, 22.69675951681078, -4.151327021478534, -3.811936822223501, -11.21234432104836, -3.387329281664682, -0.1338220704739709, -0.06083274288090479),
    @(0.2473750310259056, 0.04903106082133424, 0.1182522631433644, -0.3236351120329273, -0.09001636845753576, -0.001006874500860886, 1.164895191223389, 0.6807386034152381, 1.248617493715161, -2.410440510896091, -0.6756071513425111, -0.008203626108851017, 11.2325082463705, 4.266482512845723, 7.554762575856923, -17.70008286225281, -5.206638119442591, -0.1470323533152658, 17.3503703203209, 8.500247105664979, 9.189097835767621, -25.26364401080313, -8.982983565993445, -0.7930876848104748, -41.35448896743232, 14.373355972674, 10.27189634464736, 14.91145337324258, 2.534967181855375, -0.7371839046147971, 0.2932640248045573),
    @(0.2473750310259056, 0.04903106082133452, 0.118252263143365, -0.3236351120331354, -0.0900163684575878, -0.001006874500861462, 1.164895191232268, 0.6807386034152392, 1.248617493710722, -2.410440510897756, -0.6756071513429622, -0.008203626108856167, 11.23250824638827, 4.266482512845723, 7.554762575856925, -17.70008286226392, -5.206638119446477, -0.1470323533153526, 17.35037032039195, 8.500247105664981, 9.189097835749859, -25.26364401083865, -8.982983566006768, -0.7930876848111686, -41.3544889673968, 14.373355972674, 10.27189634463848, 14.91145337320705, 2.534967181837611, -0.7371839046170175, 0.293264024804671),
    @(0.2473750231302379, 0.04903106436627772, 0.1182522673120562, -0.3236351118629674, -0.09001636844583591, -0.001006874500740376, 1.164895165428536, 0.680738614878589, 1.248617507145291, -2.410440510097134, -0.6756071512355182, -0.008203626107611274, 11.23250816257417, 4.26648254676874, 7.554762614672881, -17.70008285343804, -5.206638117253845, -0.1470323532611953, 17.35037018230663, 8.500247147556554, 9.18909788262348, -25.26364397350888, -8.98298355466923, -0.7930876841736231, -41.35448912865589, 14.37335599316363, 10.27189637758041, 14.91145344872187, 2.534967210464306, -0.7371839009409864, 0.2932640253363843),
    @(-3.225733046709637, 1.412074165389979, 1.646627962898229, 0.1416463903768797, 0.02510842766647578, 0.0002761003780007607, -8.722718355916651, 3.609298677343222, 4.129622180468449, 0.7966373948405974, 0.1849327619434922, 0.002227341319501418, -16.115291274784, 4.808850123643972, 4.989208599260214, 4.902080451744999, 1.376955805255432, 0.03819629486497573, -9.456536531927652, 0.8474596128631194, 1.300286265304007, 5.250186924062131, 1.879489642907289, 0.1791140867390197, 22.69675947646948, -4.151327016356072, -3.811936813987968, -11.21234430216537, -3.387329274505847, -0.133822069554415, -0.06083274274800488)
)

$numRows = $newData.Count
$numCols = $newData[0].Count

# Build a true 2-D COM-compatible array from the jagged PowerShell array
$values = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $values[$r, $c] = $newData[$r][$c]
    }
}

$ws.Range("B2:AF8").Value = $values
